$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.703.87"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.599.29"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'211.30"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "1.823.39"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.600.56"
$ws.Range("E13").Value = "  -0.25%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "'65.36"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "26.677.11"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "0.0₃0759"
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").Value = "'210.21"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'7.15"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "'143.21"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'3.26"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "1.289.61"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("E36").Value = "  +1.00%  "
$ws.Range("D37").Value = "'1.50"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.0171"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("E39").Value = "  +16.62%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.784"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").Value = "'63.20"
$ws.Range("E44").Value = "  -1.24%  "
$ws.Range("D45").Value = "1.732.22"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D47").Value = "'1.58"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  -1.51%  "
